$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled city name in row 2 (F2) -- "Benaglore" -> "Bangalore"
$ws.Range("F2").Value = "Bangalore"

# Fill in the new test data row (row 3) - new user created via "add" facility
# (numeric-looking values are prefixed with a text marker so they are stored
#  as text/shared strings, matching how the rest of the row is stored, rather
#  than being auto-interpreted as numbers)
$ws.Range("B3").Value = "Tester984"
$ws.Range("C3").Value = "tester984@gmail.com"
$ws.Range("D3").Value = "'7975433984"
$ws.Range("E3").Value = "NO.10"
$ws.Range("F3").Value = "Bangalore"
$ws.Range("G3").Value = "'560023"
$ws.Range("H3").Value = "password"
$ws.Range("I3").Value = "password"

# The freshly filled row no longer carries the bordered/shaded style used
# elsewhere in the sheet -- reset it back to the workbook default.
$ws.Range("B3:I3").Style = "Normal"

# Move selection to H16 as left by the recorded session
$ws.Range("H16").Select()
